$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at the top of the data (rows 2-8), pushing the existing
# rows (previously 2..232) down to 9..239.
$ws.Rows("2:8").Insert()

# Copy the formatting (date style in col A, number style in col B) from the
# row immediately below (row 9, the old row 2) into the freshly inserted
# rows so they pick up the same cell styles / row height as the rest of the
# table.
$ws.Range("A9:B9").Copy($ws.Range("A2:B8"))
$ws.Rows("2:8").RowHeight = 18

# Populate the 7 new rows with the latest reported dates/totals (NSW second
# doses), newest first.
$ws.Range("A2").Value = 44515
$ws.Range("B2").Value = 6268391

$ws.Range("A3").Value = 44514
$ws.Range("B3").Value = 6265758

$ws.Range("A4").Value = 44513
$ws.Range("B4").Value = 6259089

$ws.Range("A5").Value = 44512
$ws.Range("B5").Value = 6243318

$ws.Range("A6").Value = 44511
$ws.Range("B6").Value = 6226830

$ws.Range("A7").Value = 44510
$ws.Range("B7").Value = 6210103

$ws.Range("A8").Value = 44509
$ws.Range("B8").Value = 6192099

# Header tweak: "second doses" -> "Second Doses"
$ws.Range("B1").Value = "Second Doses"

# Move the active selection the way the author left it.
[void]$ws.Range("D6").Select()
